$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New shared strings must land in the shared-string table in the same order
# the authored workbook has them (PR30001, Kanakambary, TC_41, TC_38), so
# the first-new-value cells are written in that order below.
$ws.Range("B7").Value = "PR30001"
$ws.Range("C7").Value = "Kanakambary"
$ws.Range("A7").Value = "TC_41"
$ws.Range("A8").Value = "TC_38"

# Remaining cells of the two new "Invalid" test case rows
$ws.Range("D7").Value = "CSA"
$ws.Range("F7").Value = "Mar_23"
$ws.Range("G7").Value = "Qualified"

$ws.Range("B8").Value = "PR30001"
$ws.Range("C8").Value = 123456
$ws.Range("D8").Value = "CSA"
$ws.Range("E8").Value = "KKEM"
$ws.Range("F8").Value = "Mar_23"
$ws.Range("G8").Value = "Qualified"

# Update the view selection state to match the authored workbook
$ws.Activate()
$ws.Range("C9").Select()
